# Applies the crypto price/volume updates described in the commit diff.
# Values are written as text (matching the inlineStr cells in the workbook);
# numeric-looking strings are apostrophe-prefixed so Excel doesn't coerce them
# to numbers, then the style is reset to Normal so no extra formatting is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.917.84"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "2.354.89"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'543.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").Value = "'134.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.82%  "

# Row 9
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").Value = "'5.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "

# Row 11
$ws.Range("E11").Value = "  -1.44%  "

# Row 12
$ws.Range("E12").Value = "  +0.27%  "

# Row 13
$ws.Range("D13").Value = "2.777.44"
$ws.Range("E13").Value = "  +0.15%  "

# Row 14
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("D15").Value = "57.858.46"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
$ws.Range("D17").Value = "2.363.34"
$ws.Range("E17").Value = "  +0.09%  "

# Row 18
$ws.Range("D18").Value = "'10.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'330.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.57%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "

# Row 21
$ws.Range("D21").Value = "'6.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "

# Row 22
$ws.Range("E22").Value = "  +0.26%  "

# Row 23
$ws.Range("D23").Value = "'62.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24
$ws.Range("E24").Value = "  -1.96%  "

# Row 25
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("D26").Value = "'8.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("E27").Value = "  -3.22%  "

# Row 28
$ws.Range("E28").Value = "  -0.30%  "

# Row 29
$ws.Range("D29").Value = "'170.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("D31").Value = "'6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("E32").Value = "  +0.01%  "

# Row 33
$ws.Range("E33").Value = "  -1.03%  "

# Row 34
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").Value = "'4.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "

# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "

# Row 37
$ws.Range("E37").Value = "  -2.81%  "

# Row 38
$ws.Range("D38").Value = "'1.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("D39").Value = "'39.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "

# Row 40
$ws.Range("D40").Value = "'142.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.78%  "

# Row 41
$ws.Range("E41").Value = "  -0.26%  "

# Row 42
$ws.Range("E42").Value = "  +0.32%  "

# Row 43
$ws.Range("D43").Value = "'289.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.55%  "

# Row 44
$ws.Range("D44").Value = "'0.0950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.04%  "

# Row 45
$ws.Range("D45").Value = "'0.0508"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "

# Row 46
$ws.Range("D46").Value = "'19.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "

# Row 47
$ws.Range("D47").Value = "'0.566"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "

# Row 48
$ws.Range("E48").Value = "  +1.62%  "

# Row 49
$ws.Range("E49").Value = "  +0.34%  "

# Row 50
$ws.Range("D50").Value = "'17.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51
$ws.Range("D51").Value = "'11.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
